$d = $word.ActiveDocument

# =====================================================================
# Change 1: paragraph 8, "{{ current_date }}" -> "{{ today() }}"
#   Rewrites the run/proofErr layout to:
#     <gramStart/> R1 "{{ "  R2 "today"  <gramEnd/>  R3 "() "  R4 "}}"
# =====================================================================
$p8 = $d.Paragraphs.Item(8)
$p8Start = $p8.Range.Start

# "current" runs from p8Start+3 .. p8Start+10, "_date" from +10 .. +15,
# " }" from +15 .. +18.  Touch across the spellStart boundary (at +3)
# with a throwaway edit so Word folds that proofErr mark away, then fix
# the word back up.
$d.Range($p8Start + 2, $p8Start + 4).Text = " t"
$d.Range($p8Start + 3, $p8Start + 10).Text = "today"

# Re-separate "{{ " from "today" into their own runs.
$d.Range($p8Start + 3, $p8Start + 8).Bold = $true
$d.Range($p8Start + 3, $p8Start + 8).Bold = $false

# Touch across the spellEnd boundary (now at +13) the same way, then
# turn "_date" into "()".
$d.Range($p8Start + 12, $p8Start + 14).Text = "X "
$d.Range($p8Start + 8, $p8Start + 13).Text = "()"

# Re-separate "() " from "}".
$d.Range($p8Start + 8, $p8Start + 11).Bold = $true
$d.Range($p8Start + 8, $p8Start + 11).Bold = $false

# =====================================================================
# Change 2: paragraph 31, merge the " ", "_____________________"
# and "              " runs into a single run.
# =====================================================================
$p31 = $d.Paragraphs.Item(31)
$p31Start = $p31.Range.Start
$mergedRange = $d.Range($p31Start, $p31Start + 36)
$mergedRange.Text = "PLACEHOLDER_TEXT_TO_FORCE_A_MERGEX00"
$mergedRange2 = $d.Range($p31Start, $p31Start + 36)
$mergedRange2.Text = " _____________________              "
